$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing A36 timestamp value (tiny float correction)
$ws.Cells.Item(36, 1).Value = 44349.85674851389

# Append new row 37 with the same date-formatted style as A36
$ws.Cells.Item(37, 1).Value = 44350.8329957329
$ws.Cells.Item(37, 1).NumberFormat = $ws.Cells.Item(36, 1).NumberFormat

$ws.Cells.Item(37, 2).Value = 75327
$ws.Cells.Item(37, 3).Value = 63469
$ws.Cells.Item(37, 4).Value = 3239
$ws.Cells.Item(37, 5).Value = 2116
$ws.Cells.Item(37, 6).Value = 1490
$ws.Cells.Item(37, 7).Value = 19755
$ws.Cells.Item(37, 8).Value = 1356
$ws.Cells.Item(37, 9).Value = 881
$ws.Cells.Item(37, 10).Value = 195
